$p = $ppt.ActivePresentation

# The "datetimeFigureOut" date placeholder (ppPlaceholderDate = 16) caches
# a literal rendering of "today" the last time the deck was saved. This
# commit simply re-caches that text (2019/1/4 -> 2020/9/25) everywhere the
# placeholder lives: the slide master, every slide layout, and the notes
# master.

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "2020/9/25"

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# Notes master.
Update-DatePlaceholder $p.NotesMaster.Shapes $newDate
